# Weekly update: prepend a new week's worth of price data (Limón, Vega
# Monumental Concepción) ahead of the existing history. Two new rows
# (1a amarillo / 2a amarillo for the new reporting date) are inserted at
# the top of the data block, pushing all later rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 259, shifting the
# existing rows 259:294 down to 261:296.
$ws.Rows("259:260").Insert()

# Seed the two new rows from the rows that used to sit at 259/260 (now at
# 261/262) so every column keeps the same labels/formatting, then
# overwrite just the cells that actually carry new-week data.
$ws.Range("A261:T261").Copy($ws.Range("A259:T259"))
$ws.Range("A262:T262").Copy($ws.Range("A260:T260"))

# New row 259: "1a amarillo"
$ws.Range("D259").Value = 44504
$ws.Range("M259").Value = 300
$ws.Range("N259").Value = 7500
$ws.Range("O259").Value = 7500
$ws.Range("P259").Value = 7500
$ws.Range("R259").Value = "Región de O'Higgins"
$ws.Range("S259").Value = 469

# New row 260: "2a amarillo"
$ws.Range("D260").Value = 44504
$ws.Range("M260").Value = 300
$ws.Range("N260").Value = 6500
$ws.Range("O260").Value = 6500
$ws.Range("P260").Value = 6500
$ws.Range("R260").Value = "Región de O'Higgins"
$ws.Range("S260").Value = 406
